# Changed syntax and processing to allow for multi-dimensional key names.
#
# The "table" placeholder syntax used a '.' to separate the table name from
# the column/key name (e.g. ${table:planData.name}). That collides with the
# '.' that would be needed to address nested / multi-dimensional keys, so
# the separator between the table name and the key name was switched to
# ':' (e.g. ${table:planData:name}).
#
# Update the three template placeholder cells on the "Plan" sheet that use
# this syntax.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

$ws.Range("B7").Value = "`${table:planData:name}"
$ws.Range("C7").Value = "`${table:planData:role}"
$ws.Range("D7").Value = "`${table:planData:days}"
